$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "fidap005.mtx"
$ws.Range("D2").Value = 0.693997343041822
$ws.Range("G2").Value = 10885.58430967119
$ws.Range("H2").Value = 3982785.614747609

$ws.Range("B3").Value = "fidap005.mtx"
$ws.Range("D3").Value = 0.8977780927497544
$ws.Range("G3").Value = 248594.0757340083
$ws.Range("H3").Value = 3491281.330207203

$ws.Range("B4").Value = "fidap005.mtx"
$ws.Range("D4").Value = 0.6855473724888353
$ws.Range("G4").Value = 0.002614121743986262
$ws.Range("H4").Value = 3983015.340353712

$ws.Range("B5").Value = "fidap005.mtx"
$ws.Range("D5").Value = 0.6855473724888353
$ws.Range("G5").Value = 1.002617541539504
$ws.Range("H5").Value = 3983015.340353712

$ws.Range("B6").Value = "fidap005.mtx"
$ws.Range("D6").Value = 0.8931493833814699
$ws.Range("G6").Value = 0.05991044653095323
$ws.Range("H6").Value = 3537273.563191011

$ws.Range("B7").Value = "fidap005.mtx"
$ws.Range("D7").Value = 0.9680559607573197
$ws.Range("F7").Value = 3669012.730238768
$ws.Range("G7").Value = 43233.30436212925
$ws.Range("H7").Value = -646.9544010491591

$ws.Range("B8").Value = "fidap005.mtx"
$ws.Range("D8").Value = 0.01293435480877391
$ws.Range("G8").Value = -76909.56717716876
$ws.Range("H8").Value = 6515955.007818108

$ws.Range("B9").Value = "fidap005.mtx"
$ws.Range("D9").Value = 0.005026103425384141
$ws.Range("G9").Value = -711314.9450195871
$ws.Range("H9").Value = 7125611.509068588

$ws.Range("B10").Value = "fidap005.mtx"
$ws.Range("D10").Value = 0.001578599958028633
$ws.Range("G10").Value = 0.002213819410285608
$ws.Range("H10").Value = 4269983.749947111

$ws.Range("B11").Value = "fidap005.mtx"
$ws.Range("D11").Value = 0.001578599958028633
$ws.Range("G11").Value = 1.002216271717798
$ws.Range("H11").Value = 4269983.749947111

$ws.Range("B12").Value = "fidap005.mtx"
$ws.Range("D12").Value = 0.009104862056452978
$ws.Range("G12").Value = 0.07888257397169395
$ws.Range("H12").Value = 3573283.453648407

$ws.Range("B13").Value = "fidap005.mtx"
$ws.Range("D13").Value = 0.02356760237874592
$ws.Range("F13").Value = 3911349.547654379
$ws.Range("G13").Value = 287735.1972457535
$ws.Range("H13").Value = -10418.42184065492

$ws.Range("B14").Value = "gr_30_30.mtx"
$ws.Range("D14").Value = 0.4751600756931221
$ws.Range("G14").Value = 0.006739676600440046
$ws.Range("H14").Value = 11.30071732523896

$ws.Range("B15").Value = "gr_30_30.mtx"
$ws.Range("D15").Value = 0.7847351381288329
$ws.Range("G15").Value = 0.3668897290933194
$ws.Range("H15").Value = 10.28027848245955

$ws.Range("B16").Value = "gr_30_30.mtx"
$ws.Range("D16").Value = 0.4589991733804813
$ws.Range("G16").Value = 0.0005912934354574122
$ws.Range("H16").Value = 11.29314559557203

$ws.Range("B17").Value = "gr_30_30.mtx"
$ws.Range("D17").Value = 0.4589991733804813
$ws.Range("G17").Value = 1.000591468283881
$ws.Range("H17").Value = 11.29314559557203

$ws.Range("B18").Value = "gr_30_30.mtx"
$ws.Range("D18").Value = 0.7684288607494869
$ws.Range("G18").Value = 0.03240812145358093
$ws.Range("H18").Value = 10.31728226696041

$ws.Range("B19").Value = "gr_30_30.mtx"
$ws.Range("D19").Value = 0.7617689471342438
$ws.Range("F19").Value = 10.83734522677404
$ws.Range("G19").Value = 0.02868986030933634
$ws.Range("H19").Value = -0.0001860185060075957

$ws.Range("B20").Value = "gr_30_30.mtx"
$ws.Range("D20").Value = 0.6140121955618729
$ws.Range("G20").Value = 0.01445828820708104
$ws.Range("H20").Value = 11.31824061810996

$ws.Range("B21").Value = "gr_30_30.mtx"
$ws.Range("D21").Value = 0.8439995573015708
$ws.Range("G21").Value = 0.3634898826266594
$ws.Range("H21").Value = 10.56654784669855

$ws.Range("B22").Value = "gr_30_30.mtx"
$ws.Range("D22").Value = 0.6028471675666579
$ws.Range("G22").Value = 0.001253504937225133
$ws.Range("H22").Value = 11.31652371167457

$ws.Range("B23").Value = "gr_30_30.mtx"
$ws.Range("D23").Value = 0.6028471675666579
$ws.Range("G23").Value = 1.001254290902909
$ws.Range("H23").Value = 11.31652371167457

$ws.Range("B24").Value = "gr_30_30.mtx"
$ws.Range("D24").Value = 0.8344806186379682
$ws.Range("G24").Value = 0.03162448067988148
$ws.Range("H24").Value = 10.59885600589429

$ws.Range("B25").Value = "gr_30_30.mtx"
$ws.Range("D25").Value = 0.8712553327746105
$ws.Range("F25").Value = 10.8664791259852
$ws.Range("G25").Value = 0.05799508875688034
$ws.Range("H25").Value = -0.0008062370472185015

$ws.Range("B26").Value = "05r0100.mtx"
$ws.Range("D26").Value = 0.9447591769470304
$ws.Range("G26").Value = -0.001339923471108601
$ws.Range("H26").Value = 12.71077434635472

$ws.Range("B27").Value = "05r0100.mtx"
$ws.Range("D27").Value = 0.9160963813348612
$ws.Range("G27").Value = -0.1052122507166789
$ws.Range("H27").Value = 13.02802165554131

$ws.Range("B28").Value = "05r0100.mtx"
$ws.Range("D28").Value = 0.9468308737415866
$ws.Range("G28").Value = -0.0001066606150576192
$ws.Range("H28").Value = 12.71131370248944

$ws.Range("B29").Value = "05r0100.mtx"
$ws.Range("D29").Value = 0.9468308737415866
$ws.Range("G29").Value = 0.9998933450729837
$ws.Range("H29").Value = 12.71131370248944

$ws.Range("B30").Value = "05r0100.mtx"
$ws.Range("D30").Value = 0.9146279331580516
$ws.Range("G30").Value = -0.008359232430789839
$ws.Range("H30").Value = 13.03546914955385

$ws.Range("B31").Value = "05r0100.mtx"
$ws.Range("D31").Value = 0.9536505262769562
$ws.Range("F31").Value = 12.73230003396526
$ws.Range("G31").Value = -0.001863170703394534
$ws.Range("H31").Value = 0.000002198517782713977

$ws.Range("B32").Value = "05r0100.mtx"
$ws.Range("D32").Value = 0.04089720730078254
$ws.Range("G32").Value = -0.002771064433890858
$ws.Range("H32").Value = 12.74844707471176

$ws.Range("B33").Value = "05r0100.mtx"
$ws.Range("D33").Value = 0.07924863152475486
$ws.Range("G33").Value = -0.06677648525432155
$ws.Range("H33").Value = 12.88482229472583

$ws.Range("B34").Value = "05r0100.mtx"
$ws.Range("D34").Value = 0.04090444282975697
$ws.Range("G34").Value = -0.0002165338169598209
$ws.Range("H34").Value = 12.74742521485896

$ws.Range("B35").Value = "05r0100.mtx"
$ws.Range("D35").Value = 0.04090444282975697
$ws.Range("G35").Value = 0.9997834896247951
$ws.Range("H35").Value = 12.74742521485896

$ws.Range("B36").Value = "05r0100.mtx"
$ws.Range("D36").Value = 0.07918161704112679
$ws.Range("G36").Value = -0.00521531502946169
$ws.Range("H36").Value = 12.88389437141706

$ws.Range("B37").Value = "05r0100.mtx"
$ws.Range("D37").Value = 0.08675933815786924
$ws.Range("F37").Value = 12.8693228114638
$ws.Range("G37").Value = -0.01726235161824638
$ws.Range("H37").Value = 0.0003450306472465552

$ws.Range("B38").Value = "bcsstm23.mtx"
$ws.Range("D38").Value = 0.5056080583310534
$ws.Range("G38").Value = 32915.90892897897
$ws.Range("H38").Value = 7842034.964505798

$ws.Range("B39").Value = "bcsstm23.mtx"
$ws.Range("D39").Value = 0.6481344316155619
$ws.Range("G39").Value = 359397.7536756572
$ws.Range("H39").Value = 7364217.656969093

$ws.Range("B40").Value = "bcsstm23.mtx"
$ws.Range("D40").Value = 0.5001900749658105
$ws.Range("G40").Value = 0.004088965653244802
$ws.Range("H40").Value = 7842022.781029087

$ws.Range("B41").Value = "bcsstm23.mtx"
$ws.Range("D41").Value = 0.5001900749658105
$ws.Range("G41").Value = 1.004097336879298
$ws.Range("H41").Value = 7842022.781029087

$ws.Range("B42").Value = "bcsstm23.mtx"
$ws.Range("D42").Value = 0.6424185998418004
$ws.Range("G42").Value = 0.04468882478301842
$ws.Range("H42").Value = 7389367.361750506

$ws.Range("B43").Value = "bcsstm23.mtx"
$ws.Range("D43").Value = 0.8448398069269917
$ws.Range("F43").Value = 6974215.319817497
$ws.Range("G43").Value = 218246.8838963124
$ws.Range("H43").Value = -8825.284522254031

$ws.Range("B44").Value = "bcsstm23.mtx"
$ws.Range("D44").Value = 0.6863846793529407
$ws.Range("G44").Value = -22871.35876629304
$ws.Range("H44").Value = 8489550.340566928

$ws.Range("B45").Value = "bcsstm23.mtx"
$ws.Range("D45").Value = 0.7826150088713799
$ws.Range("G45").Value = -186741.5519297037
$ws.Range("H45").Value = 8688712.816878704

$ws.Range("B46").Value = "bcsstm23.mtx"
$ws.Range("D46").Value = 0.6872920150900205
$ws.Range("G46").Value = -0.002743751112375878
$ws.Range("H46").Value = 8490738.227814978

$ws.Range("B47").Value = "bcsstm23.mtx"
$ws.Range("D47").Value = 0.6872920150900205
$ws.Range("G47").Value = 0.9972600095324968
$ws.Range("H47").Value = 8490738.227814978

$ws.Range("B48").Value = "bcsstm23.mtx"
$ws.Range("D48").Value = 0.7834019416537679
$ws.Range("G48").Value = -0.02239881732991944
$ws.Range("H48").Value = 8695982.461746288

$ws.Range("B49").Value = "bcsstm23.mtx"
$ws.Range("D49").Value = 0.9664987454552152
$ws.Range("F49").Value = 8995685.209173769
$ws.Range("G49").Value = -157840.6570614531
$ws.Range("H49").Value = 8435.581143447442

$ws.Range("B50").Value = "dwt_1005.mtx"
$ws.Range("D50").Value = 0.5019528886441746
$ws.Range("G50").Value = 0.01675642157169932
$ws.Range("H50").Value = 16.21183973818802

$ws.Range("B51").Value = "dwt_1005.mtx"
$ws.Range("D51").Value = 0.6280726833142906
$ws.Range("G51").Value = 0.1660399815159919
$ws.Range("H51").Value = 16.0055191530892

$ws.Range("B52").Value = "dwt_1005.mtx"
$ws.Range("D52").Value = 0.5008834592883725
$ws.Range("G52").Value = 0.001027142524775572
$ws.Range("H52").Value = 16.21191609512548

$ws.Range("B53").Value = "dwt_1005.mtx"
$ws.Range("D53").Value = 0.5008834592883725
$ws.Range("G53").Value = 1.001027670216315
$ws.Range("H53").Value = 16.21191609512548

$ws.Range("B54").Value = "dwt_1005.mtx"
$ws.Range("D54").Value = 0.6269462127925869
$ws.Range("G54").Value = 0.0101797105855257
$ws.Range("H54").Value = 16.00811241986005

$ws.Range("B55").Value = "dwt_1005.mtx"
$ws.Range("D55").Value = 0.8356592632031493
$ws.Range("F55").Value = 15.76899048563926
$ws.Range("G55").Value = 0.1193678337476338
$ws.Range("H55").Value = -0.005400600640838709

$ws.Range("B56").Value = "dwt_1005.mtx"
$ws.Range("D56").Value = 0.2092512905918695
$ws.Range("G56").Value = -0.002570307987016783
$ws.Range("H56").Value = 16.43123705900321

$ws.Range("B57").Value = "dwt_1005.mtx"
$ws.Range("D57").Value = 0.2889691199019068
$ws.Range("G57").Value = -0.02433210913379585
$ws.Range("H57").Value = 16.46051473258293

$ws.Range("B58").Value = "dwt_1005.mtx"
$ws.Range("D58").Value = 0.2091146480857755
$ws.Range("G58").Value = -0.0001564525255872695
$ws.Range("H58").Value = 16.43122087450428

$ws.Range("B59").Value = "dwt_1005.mtx"
$ws.Range("D59").Value = 0.2091146480857755
$ws.Range("G59").Value = 0.999843559712471
$ws.Range("H59").Value = 16.43122087450428

$ws.Range("B60").Value = "dwt_1005.mtx"
$ws.Range("D60").Value = 0.2888121746134026
$ws.Range("G60").Value = -0.001481156824020094
$ws.Range("H60").Value = 16.46053200549207

$ws.Range("B61").Value = "dwt_1005.mtx"
$ws.Range("D61").Value = 0.5578449002969655
$ws.Range("F61").Value = 16.5423732763387
$ws.Range("G61").Value = -0.03076904969900336
$ws.Range("H61").Value = 0.001658749512469861

$ws.Range("B62").Value = "herman4.mtx"
$ws.Range("D62").Value = 0.4504956369694701
$ws.Range("G62").Value = 0.04998131781988358
$ws.Range("H62").Value = 62.80131593208806

$ws.Range("B63").Value = "herman4.mtx"
$ws.Range("D63").Value = 0.7501546597108824
$ws.Range("G63").Value = 2.464057087186871
$ws.Range("H63").Value = 56.17071344927255

$ws.Range("B64").Value = "herman4.mtx"
$ws.Range("D64").Value = 0.4302370871312877
$ws.Range("G64").Value = 0.000789801693494947
$ws.Range("H64").Value = 62.73434082527331

$ws.Range("B65").Value = "herman4.mtx"
$ws.Range("D65").Value = 0.4302370871312877
$ws.Range("G65").Value = 1.00079011366898
$ws.Range("H65").Value = 62.73434082527331

$ws.Range("B66").Value = "herman4.mtx"
$ws.Range("D66").Value = 0.7282083816160017
$ws.Range("G66").Value = 0.03925589891238158
$ws.Range("H66").Value = 56.42624758938501

$ws.Range("B67").Value = "herman4.mtx"
$ws.Range("D67").Value = 0.7362066432967737
$ws.Range("F67").Value = 59.60769216508493
$ws.Range("G67").Value = 0.2186303516180383
$ws.Range("H67").Value = -0.001606181274268141

$ws.Range("B68").Value = "herman4.mtx"
$ws.Range("D68").Value = 0.559953221525049
$ws.Range("G68").Value = 0.1473231839705404
$ws.Range("H68").Value = 62.02677576705676

$ws.Range("B69").Value = "herman4.mtx"
$ws.Range("D69").Value = 0.7642285713989212
$ws.Range("G69").Value = 3.099431480694323
$ws.Range("H69").Value = 56.10652368570659

$ws.Range("B70").Value = "herman4.mtx"
$ws.Range("D70").Value = 0.538104210911848
$ws.Range("G70").Value = 0.002315706094106667
$ws.Range("H70").Value = 61.99549153955965

$ws.Range("B71").Value = "herman4.mtx"
$ws.Range("D71").Value = 0.538104210911848
$ws.Range("G71").Value = 1.002318389412323
$ws.Range("H71").Value = 61.99549153955965

$ws.Range("B72").Value = "herman4.mtx"
$ws.Range("D72").Value = 0.7425900839612064
$ws.Range("G72").Value = 0.04898916061623015
$ws.Range("H72").Value = 56.44146309332123

$ws.Range("B73").Value = "herman4.mtx"
$ws.Range("D73").Value = 0.8045808892920039
$ws.Range("F73").Value = 57.90391023516214
$ws.Range("G73").Value = 0.622208219057876
$ws.Range("H73").Value = -0.01079284170653033

$ws.Range("B74").Value = "dwt_1242.mtx"
$ws.Range("D74").Value = 0.7191113087771684
$ws.Range("G74").Value = 0.001761692369226648
$ws.Range("H74").Value = 9.245567265871715

$ws.Range("B75").Value = "dwt_1242.mtx"
$ws.Range("D75").Value = 0.9558644695854237
$ws.Range("G75").Value = 0.08085086586518317
$ws.Range("H75").Value = 9.035279492646803

$ws.Range("B76").Value = "dwt_1242.mtx"
$ws.Range("D76").Value = 0.7163756230604066
$ws.Range("G76").Value = 0.0001892543762774525
$ws.Range("H76").Value = 9.24552338990444

$ws.Range("B77").Value = "dwt_1242.mtx"
$ws.Range("D77").Value = 0.7163756230604066
$ws.Range("G77").Value = 1.000189272286017
$ws.Range("H77").Value = 9.24552338990444

$ws.Range("B78").Value = "dwt_1242.mtx"
$ws.Range("D78").Value = 0.9546216465232125
$ws.Range("G78").Value = 0.008696523447179704
$ws.Range("H78").Value = 9.03862732555822

$ws.Range("B79").Value = "dwt_1242.mtx"
$ws.Range("D79").Value = 0.9458481781089336
$ws.Range("F79").Value = 9.163046781593144
$ws.Range("G79").Value = 0.00593514214883235
$ws.Range("H79").Value = -0.00003794045254186948

$ws.Range("B80").Value = "dwt_1242.mtx"
$ws.Range("D80").Value = 0.7874468006916128
$ws.Range("G80").Value = 0.005119335762724632
$ws.Range("H80").Value = 9.240847873764372

$ws.Range("B81").Value = "dwt_1242.mtx"
$ws.Range("D81").Value = 0.9495657994705485
$ws.Range("G81").Value = 0.09731805223080249
$ws.Range("H81").Value = 9.064796963498605

$ws.Range("B82").Value = "dwt_1242.mtx"
$ws.Range("D82").Value = 0.785406880088627
$ws.Range("G82").Value = 0.0005491118407688352
$ws.Range("H82").Value = 9.241013549209319

$ws.Range("B83").Value = "dwt_1242.mtx"
$ws.Range("D83").Value = 0.785406880088627
$ws.Range("G83").Value = 1.000549262630275
$ws.Range("H83").Value = 9.241013549209319

$ws.Range("B84").Value = "dwt_1242.mtx"
$ws.Range("D84").Value = 0.9484605524206384
$ws.Range("G84").Value = 0.01044602274899979
$ws.Range("H84").Value = 9.067949812949317

$ws.Range("B85").Value = "dwt_1242.mtx"
$ws.Range("D85").Value = 0.9760343597285714
$ws.Range("F85").Value = 9.137649763519223
$ws.Range("G85").Value = 0.01749132614415991
$ws.Range("H85").Value = -0.0002945711995579774

$ws.Range("B86").Value = "dwt_1007.mtx"
$ws.Range("D86").Value = 0.6650730475091795
$ws.Range("G86").Value = 0.0008551751564673415
$ws.Range("H86").Value = 8.871372033654325

$ws.Range("B87").Value = "dwt_1007.mtx"
$ws.Range("D87").Value = 0.8995222168396658
$ws.Range("G87").Value = 0.02534032547217188
$ws.Range("H87").Value = 8.815320230115294

$ws.Range("B88").Value = "dwt_1007.mtx"
$ws.Range("D88").Value = 0.6641627123132186
$ws.Range("G88").Value = 0.00009620915332497082
$ws.Range("H88").Value = 8.871365573841247

$ws.Range("B89").Value = "dwt_1007.mtx"
$ws.Range("D89").Value = 0.6641627123132186
$ws.Range("G89").Value = 1.000096213781574
$ws.Range("H89").Value = 8.871365573841247

$ws.Range("B90").Value = "dwt_1007.mtx"
$ws.Range("D90").Value = 0.8988871516244913
$ws.Range("G90").Value = 0.00285179021263361
$ws.Range("H90").Value = 8.815571421198536

$ws.Range("B91").Value = "dwt_1007.mtx"
$ws.Range("D91").Value = 0.9091272063425947
$ws.Range("F91").Value = 8.842549508375454
$ws.Range("G91").Value = 0.003179761716637246
$ws.Range("H91").Value = -0.00003522100848742327

$ws.Range("B92").Value = "dwt_1007.mtx"
$ws.Range("D92").Value = 0.5663773888309865
$ws.Range("G92").Value = 0.0006702377134536849
$ws.Range("H92").Value = 8.889900737970564

$ws.Range("B93").Value = "dwt_1007.mtx"
$ws.Range("D93").Value = 0.8005175324667441
$ws.Range("G93").Value = 0.01599907089725633
$ws.Range("H93").Value = 8.857509923239158

$ws.Range("B94").Value = "dwt_1007.mtx"
$ws.Range("D94").Value = 0.5658356782872345
$ws.Range("G94").Value = 0.0000753201401171613
$ws.Range("H94").Value = 8.889894961968993

$ws.Range("B95").Value = "dwt_1007.mtx"
$ws.Range("D95").Value = 0.5658356782872345
$ws.Range("G95").Value = 1.00007532297675
$ws.Range("H95").Value = 8.889894961968993

$ws.Range("B96").Value = "dwt_1007.mtx"
$ws.Range("D96").Value = 0.800031715406112
$ws.Range("G96").Value = 0.00179826208729258
$ws.Range("H96").Value = 8.857585804786261

$ws.Range("B97").Value = "dwt_1007.mtx"
$ws.Range("D97").Value = 0.8457090467473927
$ws.Range("F97").Value = 8.868310530072936
$ws.Range("G97").Value = 0.002896032342075519
$ws.Range("H97").Value = -0.00004451589257243509

$ws.Range("B98").Value = "dwt_245.mtx"
$ws.Range("D98").Value = 0.5523584046707981
$ws.Range("G98").Value = 0.003475030716196202
$ws.Range("H98").Value = 8.396963320105652

$ws.Range("B99").Value = "dwt_245.mtx"
$ws.Range("D99").Value = 0.7987615263810328
$ws.Range("G99").Value = 0.08818713158564272
$ws.Range("H99").Value = 8.213475125809529

$ws.Range("B100").Value = "dwt_245.mtx"
$ws.Range("D100").Value = 0.5494521676266184
$ws.Range("G100").Value = 0.0004116506292178821
$ws.Range("H100").Value = 8.396742384356845

$ws.Range("B101").Value = "dwt_245.mtx"
$ws.Range("D101").Value = 0.5494521676266184
$ws.Range("G101").Value = 1.000411735368965
$ws.Range("H101").Value = 8.396742384356845

$ws.Range("B102").Value = "dwt_245.mtx"
$ws.Range("D102").Value = 0.7961148449942764
$ws.Range("G102").Value = 0.01045683324539373
$ws.Range("H102").Value = 8.215938100566452

$ws.Range("B103").Value = "dwt_245.mtx"
$ws.Range("D103").Value = 0.8606945950076306
$ws.Range("F103").Value = 8.2731977778116
$ws.Range("G103").Value = 0.0155999359686636
$ws.Range("H103").Value = -0.0002287717972163692

$ws.Range("B104").Value = "dwt_245.mtx"
$ws.Range("D104").Value = 0.3865987833025896
$ws.Range("G104").Value = 0.002364500696368224
$ws.Range("H104").Value = 8.461520916628508

$ws.Range("B105").Value = "dwt_245.mtx"
$ws.Range("D105").Value = 0.6012911007805063
$ws.Range("G105").Value = 0.0489802907705745
$ws.Range("H105").Value = 8.368256263753857

$ws.Range("B106").Value = "dwt_245.mtx"
$ws.Range("D106").Value = 0.385292808607371
$ws.Range("G106").Value = 0.0002791135380023933
$ws.Range("H106").Value = 8.461376615628321

$ws.Range("B107").Value = "dwt_245.mtx"
$ws.Range("D107").Value = 0.385292808607371
$ws.Range("G107").Value = 1.00027915249381
$ws.Range("H107").Value = 8.461376615628321

$ws.Range("B108").Value = "dwt_245.mtx"
$ws.Range("D108").Value = 0.5997391197790672
$ws.Range("G108").Value = 0.005784108157700614
$ws.Range("H108").Value = 8.36867842308002

$ws.Range("B109").Value = "dwt_245.mtx"
$ws.Range("D109").Value = 0.7126554677932113
$ws.Range("F109").Value = 8.374525695206636
$ws.Range("G109").Value = 0.01323890337410219
$ws.Range("H109").Value = -0.0002718600669433487
